$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.083099365234375
$ws.Range("B1").Value = 2.319972038269043
$ws.Range("C1").Value = 9.742735862731934
$ws.Range("D1").Value = 2.272362470626831
$ws.Range("E1").Value = 1.302870869636536
